$d = $word.ActiveDocument

# 0. Remove the pre-existing "_GoBack" bookmark that currently sits in the middle
#    of the "Njegove knjige ispod njegovih podataka" run (it will be re-added later
#    in its own new trailing paragraph), and coalesce the now-contiguous text back
#    into a single run.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$p22 = $d.Paragraphs.Item(22)
$r22 = $p22.Range
$r22.End = $r22.End - 1
$r22.Text = ""
$p22 = $d.Paragraphs.Item(22)
$r22b = $p22.Range
$r22b.End = $r22b.End - 1
$r22b.Text = "Njegove knjige ispod njegovih podataka"

# 1. Paragraph 4: "Neka piše i ime i prezime pisca, a ne samo ime"
#    -> append " unutar padajućeg izbornika"
$p4 = $d.Paragraphs.Item(4)
$r = $p4.Range
$r.End = $r.End - 1
$r.InsertAfter(" unutar padajuće")

$p4 = $d.Paragraphs.Item(4)
$r = $p4.Range
$r.End = $r.End - 1
$r.InsertAfter("g")

$p4 = $d.Paragraphs.Item(4)
$r = $p4.Range
$r.End = $r.End - 1
$r.InsertAfter(" izbornika")

# 2. Paragraph 8: "Napraviti pregled knjiga nekog žanra." -> remove trailing period
$d.Content.Find.Execute("Napraviti pregled knjiga nekog žanra.", $true, $false, $false, $false, $false, $true, 1, $false, "Napraviti pregled knjiga nekog žanra", 2)

# 3. Paragraph 22 / "Njegove knjige ispod njegovih podataka" section:
#    add new paragraphs after "Na dnu napisati koliki je ukupan broj njegovih knjiga."
#    and move the _GoBack bookmark into its own new empty trailing paragraph.

$pLast = $d.Paragraphs.Item(23)
# $pLast.Range.Text is "Na dnu napisati koliki je ukupan broj njegovih knjiga.\r"

$pLast.Range.InsertParagraphAfter()
$pNew1 = $d.Paragraphs.Item(24)
$pNew1.Range.Text = "Omogućiti kreiranje knjige na slijedeći način:"
$pNew1.Range.ListFormat.ListLevelNumber = 1

$pNew1 = $d.Paragraphs.Item(24)
$pNew1.Range.InsertParagraphAfter()
$pNew2 = $d.Paragraphs.Item(25)
$pNew2.Range.Text = "Napraviti formu na kojoj se mogu unijeti podaci o novoj knjizi – osim SlikaURL polja."
$pNew2.Range.ListFormat.ListLevelNumber = 2

$pNew2 = $d.Paragraphs.Item(25)
$pNew2.Range.InsertParagraphAfter()
$pNew3 = $d.Paragraphs.Item(26)
$pNew3.Range.Text = "Iako korisnik nije u mogućnosti unijeti taj podatak, neka se automatski unese, na serverskoj strani vrijednost: /Content/Images/<<Naslov knjige>>.jpg "
$pNew3.Range.ListFormat.ListLevelNumber = 2

# New trailing empty paragraph that will hold the relocated _GoBack bookmark.
$pNew3 = $d.Paragraphs.Item(26)
$pNew3.Range.InsertParagraphAfter()
$pNew4 = $d.Paragraphs.Item(27)
# Remove the list formatting/style inherited from the previous paragraph so it becomes a plain paragraph.
$pNew4.Range.ListFormat.RemoveNumbers()
$pNew4.set_Style("Normal")

# Bookmarks.Add misbehaves on a truly empty/collapsed range (the bookmarkEnd tag
# ends up emitted in the following paragraph), so insert a temporary placeholder
# character, anchor the bookmark around it, and then delete the placeholder again -
# the bookmark start/end stay correctly nested inside this (now empty) paragraph.
$r = $pNew4.Range
$r.End = $r.End - 1
$r.InsertAfter("X")

$pNew4 = $d.Paragraphs.Item(27)
$bmRange = $d.Range($pNew4.Range.Start, $pNew4.Range.Start + 1)
$d.Bookmarks.Add("_GoBack", $bmRange)

$pNew4 = $d.Paragraphs.Item(27)
$delRange = $d.Range($pNew4.Range.Start, $pNew4.Range.Start + 1)
$delRange.Text = ""
